$d = $word.ActiveDocument
$d.Content.Find.Execute("-1 設計類別圖", $false, $false, $false, $false, $false, $true, 1, $false, "-1 資料庫關聯表", 2)
